$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the data-organization numbers for the MESS dataset (row 22)
$ws.Range("D22").Value = 900
$ws.Range("E22").Value = 900
$ws.Range("F22").Value = 0
$ws.Range("I22").Value = 6

# Update the view to reflect where the user ended up after editing
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B29").Select()
